$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RO & CO Hearing Allocation")

$ws.Rows.Item(8).Insert()
$ws.Rows.Item(9).Copy()
$ws.Rows.Item(8).PasteSpecial(-4122)

Write-Host "done"
